# Bug fix in Eduati data files (HT115_noCTRL_meas.xlsx)
#
# Sheet1 had 43 extra stray rows (45:87) that only contained a leftover
# running index in column A, left over beyond the real data range
# (A1:N44). This trims the sheet back down to its real extent and makes
# Sheet1 the active/selected sheet (it previously was Sheet3).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")

# Remove the stray leftover rows 45:87 on Sheet1, shrinking its used
# range/dimension back down to A1:N44.
$ws1.Range("A45:A87").EntireRow.Delete()

# Sheet1 becomes the active sheet/tab (was Sheet3), with a new selection.
$ws1.Activate()
$ws1.Range("D51").Select()
